$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append at the bottom of the table (sheet lists dates most-recent-first,
# and these six rows continue the series immediately after the existing last row, 233).
$data = @(
    @("11.01.2021", 74743, 364008, 2312, 55902, 16529, 0),
    @("10.01.2021", 74665, 362432, 2300, 55502, 16863, 0),
    @("09.01.2021", 74523, 360488, 2290, 55458, 16775, 0),
    @("08.01.2021", 74233, 358330, 2278, 55277, 16678, 0),
    @("05.01.2021", 73093, 349566, 2247, 54347, 16499, 0),
    @("03.01.2021", 72573, 345504, 2196, 52979, 17398, 0)
)

$startRow = 234
$endRow = $startRow + $data.Count - 1

# Column A holds date-formatted-looking strings (e.g. "11.01.2021") that must stay
# literal text, not get auto-converted into date serial numbers when assigned.
# Force the column-A range to Text format before writing the values, then clear the
# formatting back to the sheet's default (unstyled) so the new rows look the same as
# the existing, un-styled data rows above them.
$colA = $ws.Range("A$startRow`:A$endRow")
$colA.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

$colA.ClearFormats()
